# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row => new value for column F
$updates = @{
    2  = 617
    3  = 2183
    4  = 72
    5  = 12797
    6  = 65
    7  = 113
    9  = 470
    10 = 1160
    11 = 961
    12 = 13688
    13 = 14124
    14 = 41
    15 = 170
    18 = 23
    22 = 1071
    25 = 931
    26 = 5223
    28 = 276
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
